$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header cell F1: "time_taken" label, formatted like the other header cells
$ws.Range("F1").Value = "time_taken"
$ws.Range("E1").Copy() | Out-Null
$ws.Range("F1").PasteSpecial(-4122) | Out-Null
$excel.CutCopyMode = $false

# Data cells F2:F70 with time_taken timestamps (stored as text)
$ws.Range("F2").Value = "2021-10-05 10:52:08.094750"
$ws.Range("F3").Value = "2021-10-05 10:52:08.094762"
$ws.Range("F4").Value = "2021-10-05 10:52:08.094766"
$ws.Range("F5").Value = "2021-10-05 10:52:08.094770"
$ws.Range("F6").Value = "2021-10-05 10:52:08.094773"
$ws.Range("F7").Value = "2021-10-05 10:52:08.094777"
$ws.Range("F8").Value = "2021-10-05 10:52:08.094780"
$ws.Range("F9").Value = "2021-10-05 10:52:08.094783"
$ws.Range("F10").Value = "2021-10-05 10:52:08.094786"
$ws.Range("F11").Value = "2021-10-05 10:52:08.094789"
$ws.Range("F12").Value = "2021-10-05 10:52:08.094792"
$ws.Range("F13").Value = "2021-10-05 10:52:08.094795"
$ws.Range("F14").Value = "2021-10-05 10:52:08.094798"
$ws.Range("F15").Value = "2021-10-05 10:52:08.094801"
$ws.Range("F16").Value = "2021-10-05 10:52:08.094804"
$ws.Range("F17").Value = "2021-10-05 10:52:08.094807"
$ws.Range("F18").Value = "2021-10-05 10:52:08.094811"
$ws.Range("F19").Value = "2021-10-05 10:52:08.094814"
$ws.Range("F20").Value = "2021-10-05 10:52:08.094817"
$ws.Range("F21").Value = "2021-10-05 10:52:08.094820"
$ws.Range("F22").Value = "2021-10-05 10:52:08.094823"
$ws.Range("F23").Value = "2021-10-05 10:52:08.094826"
$ws.Range("F24").Value = "2021-10-05 10:52:08.094829"
$ws.Range("F25").Value = "2021-10-05 10:52:08.094832"
$ws.Range("F26").Value = "2021-10-05 10:52:08.094835"
$ws.Range("F27").Value = "2021-10-05 10:52:08.094838"
$ws.Range("F28").Value = "2021-10-05 10:52:08.094841"
$ws.Range("F29").Value = "2021-10-05 10:52:08.094844"
$ws.Range("F30").Value = "2021-10-05 10:52:08.094847"
$ws.Range("F31").Value = "2021-10-05 10:52:08.094850"
$ws.Range("F32").Value = "2021-10-05 10:52:08.094853"
$ws.Range("F33").Value = "2021-10-05 10:52:08.094856"
$ws.Range("F34").Value = "2021-10-05 10:52:08.094859"
$ws.Range("F35").Value = "2021-10-05 10:52:08.094862"
$ws.Range("F36").Value = "2021-10-05 10:52:08.094865"
$ws.Range("F37").Value = "2021-10-05 10:52:08.094868"
$ws.Range("F38").Value = "2021-10-05 10:52:08.094871"
$ws.Range("F39").Value = "2021-10-05 10:52:08.094874"
$ws.Range("F40").Value = "2021-10-05 10:52:08.094877"
$ws.Range("F41").Value = "2021-10-05 10:52:08.094880"
$ws.Range("F42").Value = "2021-10-05 10:52:08.094884"
$ws.Range("F43").Value = "2021-10-05 10:52:08.094887"
$ws.Range("F44").Value = "2021-10-05 10:52:08.094890"
$ws.Range("F45").Value = "2021-10-05 10:52:08.094893"
$ws.Range("F46").Value = "2021-10-05 10:52:08.094896"
$ws.Range("F47").Value = "2021-10-05 10:52:08.094899"
$ws.Range("F48").Value = "2021-10-05 10:52:08.094902"
$ws.Range("F49").Value = "2021-10-05 10:52:08.094905"
$ws.Range("F50").Value = "2021-10-05 10:52:08.094908"
$ws.Range("F51").Value = "2021-10-05 10:52:08.094911"
$ws.Range("F52").Value = "2021-10-05 10:52:08.094914"
$ws.Range("F53").Value = "2021-10-05 10:52:08.094917"
$ws.Range("F54").Value = "2021-10-05 10:52:08.094920"
$ws.Range("F55").Value = "2021-10-05 10:52:08.094924"
$ws.Range("F56").Value = "2021-10-05 10:52:08.094927"
$ws.Range("F57").Value = "2021-10-05 10:52:08.094930"
$ws.Range("F58").Value = "2021-10-05 10:52:08.094933"
$ws.Range("F59").Value = "2021-10-05 10:52:08.094936"
$ws.Range("F60").Value = "2021-10-05 10:52:08.094939"
$ws.Range("F61").Value = "2021-10-05 10:52:08.094942"
$ws.Range("F62").Value = "2021-10-05 10:52:08.094945"
$ws.Range("F63").Value = "2021-10-05 10:52:08.094948"
$ws.Range("F64").Value = "2021-10-05 10:52:08.094951"
$ws.Range("F65").Value = "2021-10-05 10:52:08.094954"
$ws.Range("F66").Value = "2021-10-05 10:52:08.094958"
$ws.Range("F67").Value = "2021-10-05 10:52:08.094961"
$ws.Range("F68").Value = "2021-10-05 10:52:08.094964"
$ws.Range("F69").Value = "2021-10-05 10:52:08.094967"
$ws.Range("F70").Value = "2021-10-05 10:52:08.094970"
